$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "F"
    $ws.Cells.Item($r, 5).Value = 1
}

$ws.Range("G7").Select()
